# Apply text-cleaning (lowercasing + stemming) transformation to the
# job_title (column B) and location (column C) values, matching the
# source notebook's NLTK PorterStemmer-based preprocessing step.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jobTitles = @{
    2 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    3 = '[''nativ'', ''english'', ''teacher'', ''epik'', ''english'', ''program'', ''korea'']'
    4 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    5 = '[''peopl'', ''develop'', ''coordin'', ''ryan'']'
    6 = '[''advisori'', ''board'', ''member'', ''celal'', ''bayar'', ''univers'']'
    7 = '[''aspir'', ''human'', ''resourc'', ''specialist'']'
    8 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    9 = '[''hr'', ''senior'', ''specialist'']'
    10 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    11 = '[''seek'', ''human'', ''resourc'', ''hri'', ''generalist'', ''posit'']'
    12 = '[''student'', ''chapman'', ''univers'']'
    13 = '[''svp'', ''chro'', ''market'', ''commun'', ''csr'', ''offic'', ''engi'', ''houston'', ''woodland'', ''energi'', ''gphr'', ''sphr'']'
    14 = '[''human'', ''resourc'', ''coordin'', ''intercontinent'', ''buckhead'', ''atlanta'']'
    15 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    16 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    17 = '[''nativ'', ''english'', ''teacher'', ''epik'', ''english'', ''program'', ''korea'']'
    18 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    19 = '[''peopl'', ''develop'', ''coordin'', ''ryan'']'
    20 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    21 = '[''nativ'', ''english'', ''teacher'', ''epik'', ''english'', ''program'', ''korea'']'
    22 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    23 = '[''peopl'', ''develop'', ''coordin'', ''ryan'']'
    24 = '[''advisori'', ''board'', ''member'', ''celal'', ''bayar'', ''univers'']'
    25 = '[''aspir'', ''human'', ''resourc'', ''specialist'']'
    26 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    27 = '[''hr'', ''senior'', ''specialist'']'
    28 = '[''aspir'', ''human'', ''resourc'', ''manag'', ''student'', ''seek'', ''internship'']'
    29 = '[''seek'', ''human'', ''resourc'', ''opportun'']'
    30 = '[''aspir'', ''human'', ''resourc'', ''manag'', ''student'', ''seek'', ''internship'']'
    31 = '[''seek'', ''human'', ''resourc'', ''opportun'']'
    32 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    33 = '[''nativ'', ''english'', ''teacher'', ''epik'', ''english'', ''program'', ''korea'']'
    34 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    35 = '[''peopl'', ''develop'', ''coordin'', ''ryan'']'
    36 = '[''advisori'', ''board'', ''member'', ''celal'', ''bayar'', ''univers'']'
    37 = '[''aspir'', ''human'', ''resourc'', ''specialist'']'
    38 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    39 = '[''hr'', ''senior'', ''specialist'']'
    40 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    41 = '[''seek'', ''human'', ''resourc'', ''hri'', ''generalist'', ''posit'']'
    42 = '[''student'', ''chapman'', ''univers'']'
    43 = '[''svp'', ''chro'', ''market'', ''commun'', ''csr'', ''offic'', ''engi'', ''houston'', ''woodland'', ''energi'', ''gphr'', ''sphr'']'
    44 = '[''human'', ''resourc'', ''coordin'', ''intercontinent'', ''buckhead'', ''atlanta'']'
    45 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    46 = '[''nativ'', ''english'', ''teacher'', ''epik'', ''english'', ''program'', ''korea'']'
    47 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    48 = '[''peopl'', ''develop'', ''coordin'', ''ryan'']'
    49 = '[''advisori'', ''board'', ''member'', ''celal'', ''bayar'', ''univers'']'
    50 = '[''aspir'', ''human'', ''resourc'', ''specialist'']'
    51 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    52 = '[''hr'', ''senior'', ''specialist'']'
    53 = '[''student'', ''humber'', ''colleg'', ''aspir'', ''human'', ''resourc'', ''generalist'']'
    54 = '[''seek'', ''human'', ''resourc'', ''hri'', ''generalist'', ''posit'']'
    55 = '[''student'', ''chapman'', ''univers'']'
    56 = '[''svp'', ''chro'', ''market'', ''commun'', ''csr'', ''offic'', ''engi'', ''houston'', ''woodland'', ''energi'', ''gphr'', ''sphr'']'
    57 = '[''human'', ''resourc'', ''coordin'', ''intercontinent'', ''buckhead'', ''atlanta'']'
    58 = '[''2019'', ''ct'', ''bauer'', ''colleg'', ''busi'', ''graduat'', ''magna'', ''cum'', ''laud'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    59 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    60 = '[''peopl'', ''develop'', ''coordin'', ''ryan'']'
    61 = '[''aspir'', ''human'', ''resourc'', ''specialist'']'
    62 = '[''hr'', ''senior'', ''specialist'']'
    63 = '[''seek'', ''human'', ''resourc'', ''hri'', ''generalist'', ''posit'']'
    64 = '[''student'', ''chapman'', ''univers'']'
    65 = '[''svp'', ''chro'', ''market'', ''commun'', ''csr'', ''offic'', ''engi'', ''houston'', ''woodland'', ''energi'', ''gphr'', ''sphr'']'
    66 = '[''human'', ''resourc'', ''coordin'', ''intercontinent'', ''buckhead'', ''atlanta'']'
    67 = '[''experienc'', ''retail'', ''manag'', ''aspir'', ''human'', ''resourc'', ''profession'']'
    68 = '[''human'', ''resourc'', ''staf'', ''recruit'', ''profession'']'
    69 = '[''human'', ''resourc'', ''specialist'', ''luxottica'']'
    70 = '[''director'', ''human'', ''resourc'', ''north'', ''america'', ''group'', ''beneteau'']'
    71 = '[''retir'', ''armi'', ''nation'', ''guard'', ''recruit'', ''offic'', ''manag'', ''seek'', ''posit'', ''human'', ''resourc'']'
    72 = '[''human'', ''resourc'', ''generalist'', ''scottmadden'', ''inc'']'
    73 = '[''busi'', ''manag'', ''major'', ''aspir'', ''human'', ''resourc'', ''manag'']'
    74 = '[''aspir'', ''human'', ''resourc'', ''manag'', ''seek'', ''internship'', ''human'', ''resourc'']'
    75 = '[''human'', ''resourc'', ''profession'']'
    76 = '[''nortia'', ''staf'', ''seek'', ''human'', ''resourc'', ''payrol'', ''administr'', ''profession'', ''408'', ''7092621'']'
    77 = '[''aspir'', ''human'', ''resourc'', ''profession'', ''passion'', ''help'', ''creat'', ''inclus'', ''engag'', ''work'', ''environ'']'
    78 = '[''human'', ''resourc'', ''conflict'', ''manag'', ''polici'', ''procedurestal'', ''managementbenefit'', ''compens'']'
    79 = '[''human'', ''resourc'', ''generalist'', ''schwan'']'
    80 = '[''liber'', ''art'', ''major'', ''aspir'', ''human'', ''resourc'', ''analyst'']'
    81 = '[''junior'', ''me'', ''engin'', ''inform'', ''system'']'
    82 = '[''senior'', ''human'', ''resourc'', ''busi'', ''partner'', ''heil'', ''environment'']'
    83 = '[''aspir'', ''human'', ''resourc'', ''profession'', ''energet'', ''teamfocus'', ''leader'']'
    84 = '[''hr'', ''manag'', ''endemol'', ''shine'', ''north'', ''america'']'
    85 = '[''human'', ''resourc'', ''profession'', ''world'', ''leader'', ''gi'', ''softwar'']'
    86 = '[''rrp'', ''brand'', ''portfolio'', ''execut'', ''jti'', ''japan'', ''tobacco'', ''intern'']'
    87 = '[''inform'', ''system'', ''specialist'', ''programm'', ''love'', ''data'', ''organ'']'
    88 = '[''bachelor'', ''scienc'', ''biolog'', ''victoria'', ''univers'', ''wellington'']'
    89 = '[''human'', ''resourc'', ''manag'', ''major'']'
    90 = '[''director'', ''human'', ''resourc'', ''ey'']'
    91 = '[''undergradu'', ''research'', ''assist'', ''styczynski'', ''lab'']'
    92 = '[''lead'', ''offici'', ''western'', ''illinoi'', ''univers'']'
    93 = '[''seek'', ''employ'', ''opportun'', ''within'', ''custom'', ''servic'', ''patient'', ''care'']'
    94 = '[''admiss'', ''repres'', ''commun'', ''medic'', ''center'', ''long'', ''beach'']'
    95 = '[''seek'', ''human'', ''resourc'', ''opportun'', ''open'', ''travel'', ''reloc'']'
    96 = '[''student'', ''westfield'', ''state'', ''univers'']'
    97 = '[''student'', ''indiana'', ''univers'', ''kokomo'', ''busi'', ''manag'', ''retail'', ''manag'', ''delphi'', ''hardwar'', ''paint'']'
    98 = '[''aspir'', ''human'', ''resourc'', ''profession'']'
    99 = '[''student'']'
    100 = '[''seek'', ''human'', ''resourc'', ''posit'']'
    101 = '[''aspir'', ''human'', ''resourc'', ''manag'', ''graduat'', ''may'', ''2020'', ''seek'', ''entrylevel'', ''human'', ''resourc'', ''posit'', ''st'', ''loui'']'
    102 = '[''human'', ''resourc'', ''generalist'', ''loparex'']'
    103 = '[''busi'', ''intellig'', ''analyt'', ''travel'']'
    104 = '[''alway'', ''set'', ''success'']'
    105 = '[''director'', ''administr'', ''excel'', ''log'']'
}

$locations = @{
    2 = '[''houston'', ''texa'']'
    3 = '[''kanada'']'
    4 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    5 = '[''denton'', ''texa'']'
    6 = '[''i̇zmir'', ''türkiy'']'
    7 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    8 = '[''kanada'']'
    9 = '[''san'', ''francisco'', ''bay'', ''area'']'
    10 = '[''kanada'']'
    11 = '[''greater'', ''philadelphia'', ''area'']'
    12 = '[''lake'', ''forest'', ''california'']'
    13 = '[''houston'', ''texa'', ''area'']'
    14 = '[''atlanta'', ''georgia'']'
    15 = '[''houston'', ''texa'']'
    16 = '[''houston'', ''texa'']'
    17 = '[''kanada'']'
    18 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    19 = '[''denton'', ''texa'']'
    20 = '[''houston'', ''texa'']'
    21 = '[''kanada'']'
    22 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    23 = '[''denton'', ''texa'']'
    24 = '[''i̇zmir'', ''türkiy'']'
    25 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    26 = '[''kanada'']'
    27 = '[''san'', ''francisco'', ''bay'', ''area'']'
    28 = '[''houston'', ''texa'', ''area'']'
    29 = '[''chicago'', ''illinoi'']'
    30 = '[''houston'', ''texa'', ''area'']'
    31 = '[''chicago'', ''illinoi'']'
    32 = '[''houston'', ''texa'']'
    33 = '[''kanada'']'
    34 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    35 = '[''denton'', ''texa'']'
    36 = '[''i̇zmir'', ''türkiy'']'
    37 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    38 = '[''kanada'']'
    39 = '[''san'', ''francisco'', ''bay'', ''area'']'
    40 = '[''kanada'']'
    41 = '[''greater'', ''philadelphia'', ''area'']'
    42 = '[''lake'', ''forest'', ''california'']'
    43 = '[''houston'', ''texa'', ''area'']'
    44 = '[''atlanta'', ''georgia'']'
    45 = '[''houston'', ''texa'']'
    46 = '[''kanada'']'
    47 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    48 = '[''denton'', ''texa'']'
    49 = '[''i̇zmir'', ''türkiy'']'
    50 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    51 = '[''kanada'']'
    52 = '[''san'', ''francisco'', ''bay'', ''area'']'
    53 = '[''kanada'']'
    54 = '[''greater'', ''philadelphia'', ''area'']'
    55 = '[''lake'', ''forest'', ''california'']'
    56 = '[''houston'', ''texa'', ''area'']'
    57 = '[''atlanta'', ''georgia'']'
    58 = '[''houston'', ''texa'']'
    59 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    60 = '[''denton'', ''texa'']'
    61 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    62 = '[''san'', ''francisco'', ''bay'', ''area'']'
    63 = '[''greater'', ''philadelphia'', ''area'']'
    64 = '[''lake'', ''forest'', ''california'']'
    65 = '[''houston'', ''texa'', ''area'']'
    66 = '[''atlanta'', ''georgia'']'
    67 = '[''austin'', ''texa'', ''area'']'
    68 = '[''jackson'', ''mississippi'', ''area'']'
    69 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    70 = '[''greater'', ''grand'', ''rapid'', ''michigan'', ''area'']'
    71 = '[''virginia'', ''beach'', ''virginia'']'
    72 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    73 = '[''monro'', ''louisiana'', ''area'']'
    74 = '[''houston'', ''texa'', ''area'']'
    75 = '[''greater'', ''boston'', ''area'']'
    76 = '[''san'', ''jose'', ''california'']'
    77 = '[''new'', ''york'', ''new'', ''york'']'
    78 = '[''dallasfort'', ''worth'', ''area'']'
    79 = '[''amerika'', ''birleşik'', ''devletleri'']'
    80 = '[''baton'', ''roug'', ''louisiana'', ''area'']'
    81 = '[''myrtl'', ''beach'', ''south'', ''carolina'', ''area'']'
    82 = '[''chattanooga'', ''tennesse'', ''area'']'
    83 = '[''austin'', ''texa'', ''area'']'
    84 = '[''lo'', ''angel'', ''california'']'
    85 = '[''highland'', ''california'']'
    86 = '[''greater'', ''philadelphia'', ''area'']'
    87 = '[''gaithersburg'', ''maryland'']'
    88 = '[''baltimor'', ''maryland'']'
    89 = '[''milpita'', ''california'']'
    90 = '[''greater'', ''atlanta'', ''area'']'
    91 = '[''greater'', ''atlanta'', ''area'']'
    92 = '[''greater'', ''chicago'', ''area'']'
    93 = '[''torranc'', ''california'']'
    94 = '[''long'', ''beach'', ''california'']'
    95 = '[''amerika'', ''birleşik'', ''devletleri'']'
    96 = '[''bridgewat'', ''massachusett'']'
    97 = '[''lafayett'', ''indiana'']'
    98 = '[''kokomo'', ''indiana'', ''area'']'
    99 = '[''houston'', ''texa'', ''area'']'
    100 = '[''la'', ''vega'', ''nevada'', ''area'']'
    101 = '[''cape'', ''girardeau'', ''missouri'']'
    102 = '[''raleighdurham'', ''north'', ''carolina'', ''area'']'
    103 = '[''greater'', ''new'', ''york'', ''citi'', ''area'']'
    104 = '[''greater'', ''lo'', ''angel'', ''area'']'
    105 = '[''kati'', ''texa'']'
}

foreach ($r in $jobTitles.Keys) {
    $ws.Cells.Item($r, 2).Value = $jobTitles[$r]
    $ws.Cells.Item($r, 3).Value = $locations[$r]
}
